# Auto-generated script to apply cryptos.xlsx price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '69.354.18'
$ws.Range('E2').Value = '  +1.94%  '

# Row 3
$ws.Range('D3').Value = '3.781.21'
$ws.Range('E3').Value = '  +0.01%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  +0.37%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '622.70'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +3.91%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.68'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.51%  '

# Row 7
$ws.Range('D7').Value = '3.778.81'
$ws.Range('E7').Value = '  -0.01%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '1.00'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -0.14%  '

# Row 10
$ws.Range('E10').Value = '  +2.42%  '

# Row 11
$ws.Range('E11').Value = '  +1.31%  '

# Row 12
$ws.Range('E12').Value = '  +0.87%  '

# Row 13
$ws.Range('E13').Value = '  +0.45%  '

# Row 14
$ws.Range('E14').Value = '  +0.92%  '

# Row 15
$ws.Range('D15').Value = '4.413.59'
$ws.Range('E15').Value = '  -0.02%  '

# Row 16
$ws.Range('D16').Value = '3.767.07'
$ws.Range('E16').Value = '  -1.07%  '

# Row 17
$ws.Range('D17').Value = '69.315.00'
$ws.Range('E17').Value = '  +2.00%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '17.67'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -3.17%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.11'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +1.43%  '

# Row 20
$ws.Range('E20').Value = '  -1.04%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '468.28'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +2.40%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.64'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.51%  '

# Row 23
$ws.Range('E23').Value = '  +1.19%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.0000150'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +5.40%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.39'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.77%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.02'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.73%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.16'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +3.86%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.04'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +1.44%  '

# Row 29
$ws.Range('E29').Value = '  +0.01%  '

# Row 30
$ws.Range('D30').Value = '3.930.03'
$ws.Range('E30').Value = '  +0.03%  '

# Row 31
$ws.Range('E31').Value = '  +3.81%  '

# Row 32
$ws.Range('E32').Value = '  +1.61%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '7.32'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +1.35%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '28.86'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.41%  '

# Row 35
$ws.Range('E35').Value = '  -0.05%  '

# Row 36
$ws.Range('D36').Value = '3.731.17'
$ws.Range('E36').Value = '  +0.08%  '

# Row 37
$ws.Range('E37').Value = '  +0.74%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.162'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +14.36%  '

# Row 39
$ws.Range('E39').Value = '  +2.96%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.40'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +7.28%  '

# Row 41
$ws.Range('E41').Value = '  +0.56%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.969'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -1.10%  '

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.999'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.01%  '

# Row 45
$ws.Range('B45').Value = 'TheGraph'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.300'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +1.77%  '

# Row 46
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '154.34'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +1.09%  '

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.26'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -0.10%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '46.82'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.97%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.91'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.55%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '8.43'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +1.69%  '

# Row 51
$ws.Range('E51').Value = '  +0.12%  '
